# Add a new "canonical SMILES" column (D) next to the existing
# "canonical isomeric SMILES" column (C), populated with the isomeric
# SMILES stripped of the stereo-bond markers ("/" and "\").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in D2 (the new cell automatically picks up the row's
# existing formatting, same as the neighboring header cells).
$ws.Cells.Item(2, 4).Value = "canonical SMILES"

# Give column D a sensible width, similar to the other data columns
# (matches the 38.42578125 width used in the canonical workbook).
$ws.Columns.Item(4).ColumnWidth = 37.67

# Data rows are 3..33 (molecule rows); B = microstate ID, C = canonical
# isomeric SMILES. Fill D with the canonical (non-isomeric) SMILES.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
for ($r = 3; $r -le $lastRow; $r++) {
    $iso = $ws.Cells.Item($r, 3).Value()
    if ($iso -ne $null) {
        $canon = $iso.Replace("/", "").Replace("\", "")
        $ws.Cells.Item($r, 4).Value = $canon
    }
}
